$x = $ppt.ActivePresentation.Slides.Item(999)
Write-Host ("Type: " + $x.GetType())
Write-Host ("Value: " + $x)
